$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Prefix the in-progress / temporary sample names with "TMP_" so they are
# clearly flagged as not-yet-finalized entries in the readme sheet.
$ws1.Range("A9").Value  = "TMP_FW_SOURCE_HR4"
$ws1.Range("A10").Value = "TMP_ESTUARY_BARGE_HR8"
$ws1.Range("A11").Value = "TMP_FW_SOURCE_HR5"
$ws1.Range("A12").Value = "TMP_FW_SOURCE_HR7"
$ws1.Range("A13").Value = "TMP_SW_SOURCE_HR7"
$ws1.Range("A14").Value = "TMP_FW_SOURCE_HR0"

# Bring Sheet1 to the front (it was "Dilution sheet" before) and leave the
# selection where the author last clicked.
$null = $ws1.Activate()
$null = $ws1.Range("A18").Select()
